$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# pseudo-F (column E) updates for rows 2-21
$ws.Range("E2").Value = 4.22099180635082
$ws.Range("E3").Value = 7.144241145310827
$ws.Range("E4").Value = 10.13108699498125
$ws.Range("E5").Value = 6.044788046584133
$ws.Range("E6").Value = 3.581301113084161
$ws.Range("E7").Value = 10.8476464422071
$ws.Range("E8").Value = 5.033496363401198
$ws.Range("E9").Value = 2.881999722414614
$ws.Range("E10").Value = 2.380743607023501
$ws.Range("E11").Value = 1.507498411901389
$ws.Range("E12").Value = 16.29361754972921
$ws.Range("E13").Value = 15.65402537301194
$ws.Range("E14").Value = 9.560102567112489
$ws.Range("E15").Value = 4.14974846932162
$ws.Range("E16").Value = 7.322697083789895
$ws.Range("E17").Value = 6.251493246914585
$ws.Range("E18").Value = 7.588921594740497
$ws.Range("E19").Value = 7.70977198010405
$ws.Range("E20").Value = 6.923125790138299
$ws.Range("E21").Value = 4.242178731174844

# p-value (column F) and q-value (column G) updates for rows 9-11
$ws.Range("F9").Value = 0.002
$ws.Range("G9").Value = 0.0025

$ws.Range("F10").Value = 0.014
$ws.Range("G10").Value = 0.01555555555555556

$ws.Range("F11").Value = 0.133
$ws.Range("G11").Value = 0.133

Write-Host "done"
